$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in previously-blank TotalConfirmedNewCases (G) / TotalNewDeaths (I)
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 0

$ws.Range("G7").Value = 5
$ws.Range("I7").Value = 0

# ---------------------------------------------------------------------------
# 2. Rename header L1
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

# ---------------------------------------------------------------------------
# 3. Append new row 8
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 71
$ws.Range("B8").Value = 119
$ws.Range("C8").Value = "NORTHERN AFRICA                    "
$ws.Range("D8").Value = Get-Date -Year 2020 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Range("E8").Value = "Libya"
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "Local transmission"
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 5331

# ---------------------------------------------------------------------------
# 4. Column widths: A:O = 27
# ---------------------------------------------------------------------------
$ws.Range("A1:O8").ColumnWidth = 27

# ---------------------------------------------------------------------------
# 5. Formatting: center horizontal/vertical alignment everywhere (A1:O8),
#    and a dedicated date format + alignment for the Date column (D1:D8).
#    Build each combo once on a scratch cell, then PasteSpecial-Formats it
#    across the destination range so only a single new style is minted.
# ---------------------------------------------------------------------------
$ws.Range("Z1").HorizontalAlignment = -4108
$ws.Range("Z1").VerticalAlignment = -4108

$ws.Range("Z2").NumberFormat = "yyyy-mm-dd;"
$ws.Range("Z2").HorizontalAlignment = -4108
$ws.Range("Z2").VerticalAlignment = -4108

$ws.Range("Z1").Copy()
$ws.Range("A1:O8").PasteSpecial(-4122)

$ws.Range("Z2").Copy()
$ws.Range("D1:D8").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()

$excel.CutCopyMode = 0
$ws.Range("A1").Select()
